$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.04312
$ws.Range("H2").Value2 = 0.12936
$ws.Range("I2").Value2 = 0.06332032271043876
$ws.Range("J2").Value2 = 0.06332032271043876
$ws.Range("M2").Value2 = 0.4724393333333333
$ws.Range("N2").Value2 = 1.417318
$ws.Range("O2").Value2 = 0.04281008045734225
$ws.Range("P2").Value2 = 0.04281008045734226
$ws.Range("Q2").Value2 = 0.02037158405333333
$ws.Range("R2").Value2 = 0.18334425648
$ws.Range("S2").Value2 = 0.002710748109818759
$ws.Range("T2").Value2 = 0.00271074810981876
# Row 3
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.04312
$ws.Range("H3").Value2 = 0.12936
$ws.Range("I3").Value2 = 0.06332032271043876
$ws.Range("J3").Value2 = 0.06332032271043876
$ws.Range("O3").Value2 = 0.6604099741840944
$ws.Range("P3").Value2 = 0.6604099741840945
$ws.Range("Q3").Value2 = 0.3142623689333333
$ws.Range("R3").Value2 = 2.8283613204
$ws.Range("S3").Value2 = 0.04181737268652939
$ws.Range("T3").Value2 = 0.04181737268652939
# Row 4
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.04312
$ws.Range("H4").Value2 = 0.12936
$ws.Range("I4").Value2 = 0.06332032271043876
$ws.Range("J4").Value2 = 0.06332032271043876
$ws.Range("O4").Value2 = 0.2967799453585633
$ws.Range("P4").Value2 = 0.2967799453585634
$ws.Range("Q4").Value2 = 0.1412255603733333
$ws.Range("R4").Value2 = 1.27103004336
$ws.Range("S4").Value2 = 0.01879220191409061
$ws.Range("T4").Value2 = 0.01879220191409062
# Row 5
$ws.Range("I5").Value2 = 0.5780441577995699
$ws.Range("J5").Value2 = 0.5780441577995699
$ws.Range("M5").Value2 = 0.4724393333333333
$ws.Range("N5").Value2 = 1.417318
$ws.Range("O5").Value2 = 0.04281008045734225
$ws.Range("P5").Value2 = 0.04281008045734226
$ws.Range("Q5").Value2 = 0.1859699168148888
$ws.Range("R5").Value2 = 1.673729251334
$ws.Range("S5").Value2 = 0.02474611690329623
$ws.Range("T5").Value2 = 0.02474611690329623
# Row 6
$ws.Range("I6").Value2 = 0.5780441577995699
$ws.Range("J6").Value2 = 0.5780441577995699
$ws.Range("O6").Value2 = 0.6604099741840944
$ws.Range("P6").Value2 = 0.6604099741840945
$ws.Range("S6").Value2 = 0.3817461273296806
$ws.Range("T6").Value2 = 0.3817461273296806
# Row 7
$ws.Range("I7").Value2 = 0.5780441577995699
$ws.Range("J7").Value2 = 0.5780441577995699
$ws.Range("O7").Value2 = 0.2967799453585633
$ws.Range("P7").Value2 = 0.2967799453585634
$ws.Range("S7").Value2 = 0.1715519135665931
$ws.Range("T7").Value2 = 0.1715519135665931
# Row 8
$ws.Range("H8").Value2 = 0.7326729999999999
$ws.Range("I8").Value2 = 0.3586355194899915
$ws.Range("J8").Value2 = 0.3586355194899914
$ws.Range("M8").Value2 = 0.4724393333333333
$ws.Range("N8").Value2 = 1.417318
$ws.Range("O8").Value2 = 0.04281008045734225
$ws.Range("P8").Value2 = 0.04281008045734226
$ws.Range("Q8").Value2 = 0.1153811812237777
$ws.Range("R8").Value2 = 1.038430631014
$ws.Range("S8").Value2 = 0.01535321544422727
$ws.Range("T8").Value2 = 0.01535321544422727
# Row 9
$ws.Range("H9").Value2 = 0.7326729999999999
$ws.Range("I9").Value2 = 0.3586355194899915
$ws.Range("J9").Value2 = 0.3586355194899914
$ws.Range("O9").Value2 = 0.6604099741840944
$ws.Range("P9").Value2 = 0.6604099741840945
$ws.Range("Q9").Value2 = 1.779928514482777
$ws.Range("S9").Value2 = 0.2368464741678845
$ws.Range("T9").Value2 = 0.2368464741678845
# Row 10
$ws.Range("H10").Value2 = 0.7326729999999999
$ws.Range("I10").Value2 = 0.3586355194899915
$ws.Range("J10").Value2 = 0.3586355194899914
$ws.Range("O10").Value2 = 0.2967799453585633
$ws.Range("P10").Value2 = 0.2967799453585634
$ws.Range("Q10").Value2 = 0.799877512333111
$ws.Range("R10").Value2 = 7.198897610997999
$ws.Range("S10").Value2 = 0.1064358298778796
$ws.Range("T10").Value2 = 0.1064358298778796
